# Generate Report for Handback
# Update the localization-status report: the de-de (and zh-cn) handback is now
# in sync with en-US, so refresh status + handback timestamps and clear the
# stale "handback file is not latest" error detail.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$status = "Handed back: in sync with en-US"

# Overview sheet: both locale status columns move to "handed back" state.
$wsOverview.Range("E2").Value = $status
$wsOverview.Range("F2").Value = $status

# zh-cn detail sheet: status, refreshed handback datetime, cleared error detail.
$wsZhCn.Range("C2").Value = $status
$wsZhCn.Range("K2").Value = "2016-09-04 10:54:28"
$wsZhCn.Range("P2").Value = ""

# de-de detail sheet: status, refreshed handback datetime, cleared error detail.
$wsDeDe.Range("C2").Value = $status
$wsDeDe.Range("K2").Value = "2016-09-04 10:54:35"
$wsDeDe.Range("P2").Value = ""

# Widen the Status columns (now holding longer text) and shrink the now-empty
# Error Detail columns to match their content.
$wsOverview.Columns.Item(5).ColumnWidth = 29.16666666667
$wsOverview.Columns.Item(6).ColumnWidth = 29.16666666667

$wsZhCn.Columns.Item(3).ColumnWidth = 29.16666666667
$wsZhCn.Columns.Item(16).ColumnWidth = 12.83333333333

$wsDeDe.Columns.Item(3).ColumnWidth = 29.16666666667
$wsDeDe.Columns.Item(16).ColumnWidth = 12.83333333333
